$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.354.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.626.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9995'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '302.44'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3761'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3625'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.42'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08154'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.220'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.24'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.465'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.87%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001238'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.37%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.293'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.606.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.32'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06928'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.59%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.543'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.348.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.475'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.063'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.12'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.14'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.276'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.80'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.794.03'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.703'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.139'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.053'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +10.91%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.19'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02763'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08765'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2489'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.07121'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.983'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.49%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.332'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.84'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.03'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6449'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9988'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.274'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.958'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07969'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.56%  '
